$d = $word.ActiveDocument

# Walk every story (main text, headers, footers, ...) looking for the
# inline pictures whose internal drawing name needs to be swapped:
#   - the two Pearson logo pictures in the footers: image2.png -> image1.png
#   - the two BTEC logo pictures in the headers:   image1.jpg -> image2.jpg
# (WdStoryType: 7/10 = primary/first-page header, 9/11 = primary/first-page footer)
foreach ($story in $d.StoryRanges) {
    if ($story.InlineShapes.Count -gt 0) {
        $shape = $story.InlineShapes.Item(1)

        # Re-fetch the shape through Selection before renaming it - footer-anchored
        # inline shapes need a "live" selection handle for the Name property to apply.
        $null = $shape.Select()
        $target = $word.Selection.InlineShapes.Item(1)

        if ($story.StoryType -eq 9 -or $story.StoryType -eq 11) {
            $target.Name = "image1.png"
        } elseif ($story.StoryType -eq 7 -or $story.StoryType -eq 10) {
            $target.Name = "image2.jpg"
        }
    }
}
